$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.927.79"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "3.120.75"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'526.58"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").Value = "'141.47"
$ws.Range("E6").Value = "  +1.16%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.125.27"
$ws.Range("E8").Value = "  +2.43%  "
$ws.Range("D9").Value = "'0.434"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").Value = "'7.28"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("E12").Value = "  +3.42%  "
$ws.Range("D13").Value = "3.650.95"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").Value = "'26.41"
$ws.Range("E15").Value = "  +4.23%  "
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").Value = "58.049.91"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "3.112.44"
$ws.Range("E18").Value = "  +2.02%  "
$ws.Range("D19").Value = "'6.15"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "'12.90"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "'8.11"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "'337.85"
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +2.64%  "
$ws.Range("D25").Value = "'66.72"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").Value = "'0.169"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  +3.64%  "
$ws.Range("D29").Value = "'6.60"
$ws.Range("E29").Value = "  +4.53%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'7.27"
$ws.Range("E31").Value = "  +1.49%  "
$ws.Range("E32").Value = "  +3.22%  "
$ws.Range("E33").Value = "  +3.47%  "
$ws.Range("D34").Value = "'21.03"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("D35").Value = "'154.16"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("E36").Value = "  +5.49%  "
$ws.Range("D37").Value = "'6.13"
$ws.Range("E37").Value = "  +3.91%  "
$ws.Range("D38").Value = "'26.96"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").Value = "'1.32"
$ws.Range("E39").Value = "  +3.08%  "
$ws.Range("D40").Value = "'0.0670"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").Value = "3.157.93"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("D42").Value = "'0.689"
$ws.Range("E42").Value = "  +5.74%  "
# Row 43: reorder -> Stacks
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.52"
$ws.Range("E43").Value = "  +10.67%  "
# Row 44: reorder -> Filecoin
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'3.92"
$ws.Range("E44").Value = "  +0.53%  "
# Row 45: reorder -> OKB
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'37.07"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "2.302.68"
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("E49").Value = "  +8.04%  "
$ws.Range("D50").Value = "'21.05"
$ws.Range("E50").Value = "  +4.27%  "
$ws.Range("D51").Value = "'6.03"
$ws.Range("E51").Value = "  +2.99%  "
